# Update the "展览" (sheet 1) and "全部类型" (sheet 4) tables:
#  - remove the oldest entry (2024-06-01, 丽水·动漫游戏展), shifting the
#    remaining three rows up by one
#  - bump the "想去人数" (F column) counter by 1 on the two rows whose
#    underlying event data changed position
#  - renumber the index column (A) to 1..3

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the first data row (2024-06-01 / 丽水·动漫游戏展); remaining
    # rows shift up automatically.
    $ws.Rows.Item(2).Delete()

    # Re-index the "序号" column now that rows have shifted up.
    $ws.Range("A2").Value = 1
    $ws.Range("A3").Value = 2
    $ws.Range("A4").Value = 3

    # Attendance counts ("想去人数") increased by 1 for these two events.
    $ws.Range("F2").Value = 22
    $ws.Range("F3").Value = 46
}
